$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.66"
$ws.Range("E2").Value = "'5.16%"

$ws.Range("D3").Value = "'27.36"
$ws.Range("E3").Value = "'-3.40%"

$ws.Range("E4").Value = "'-1.33%"

$ws.Range("D5").Value = "'0.05931"
$ws.Range("E5").Value = "'3.88%"

$ws.Range("D6").Value = "'6.708"
$ws.Range("E6").Value = "'1.01%"

$ws.Range("D7").Value = "'0.8701"
$ws.Range("E7").Value = "'0.82%"

$ws.Range("D8").Value = "'0.9983"
$ws.Range("E8").Value = "'12.90%"

$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01051"
$ws.Range("E9").Value = "'1,660.92%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("E10").Value = "'2.05%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07171"
$ws.Range("E11").Value = "'1.33%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03148"
$ws.Range("E12").Value = "'0.01%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09231"
$ws.Range("E13").Value = "'0.05%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001551"
$ws.Range("E14").Value = "'1.79%"

$ws.Range("D15").Value = "'0.005851"
$ws.Range("E15").Value = "'-3.50%"

$ws.Range("D16").Value = "'3.502"
$ws.Range("E16").Value = "'0.17%"

$ws.Range("E17").Value = "'1.66%"

$ws.Range("E18").Value = "'2.50%"

$ws.Range("E19").Value = "'-0.71%"

$ws.Range("D20").Value = "'0.03557"
$ws.Range("E20").Value = "'6.72%"

$ws.Range("D21").Value = "'0.1305"
$ws.Range("E21").Value = "'-0.31%"

$ws.Range("D22").Value = "'3.520"
$ws.Range("E22").Value = "'1.01%"

$ws.Range("E23").Value = "'2.26%"

$ws.Range("E24").Value = "'1.43%"

$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'-0.41%"

$ws.Range("E26").Value = "'8.65%"

$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.06%"

$ws.Range("D28").Value = "'0.0001936"
$ws.Range("E28").Value = "'34.03%"

$ws.Range("D40").Value = "'0.03835"
$ws.Range("E40").Value = "'0.55%"

$ws.Range("D41").Value = "'0.006574"
$ws.Range("E41").Value = "'14.57%"

$ws.Range("D42").Value = "'0.1106"
$ws.Range("E42").Value = "'3.46%"

$ws.Range("D43").Value = "'0.002277"
$ws.Range("E43").Value = "'3.55%"

$ws.Range("D44").Value = "'0.01071"
$ws.Range("E44").Value = "'12.77%"

$ws.Range("D45").Value = "'0.00005444"
$ws.Range("E45").Value = "'3.30%"

$ws.Range("E46").Value = "'-0.05%"

$ws.Range("D47").Value = "'0.1090"
$ws.Range("E47").Value = "'22.38%"

$ws.Range("D48").Value = "'0.002233"
$ws.Range("E48").Value = "'-0.99%"

$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.05%"

$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.05%"
